$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Employee number" / "Employee Number" header columns to
# "Departement Simple" and the "Text (7 digits)" type hint to "Tekst".
$ws.Range("F8").Value = "Tekst"
$ws.Range("F7").Value = "Departement Simple"

# Clear the sample employee number value in F9, keeping its formatting.
$ws.Range("F9").ClearContents()

# Update the second header block (row 13) to match the renamed column,
# and apply the same text number format used by the other F column cells.
$ws.Range("F13").Value = "Departement Simple"
$ws.Range("F13").NumberFormat = "@"

# Widen column F slightly to fit the new, longer header text.
# (17.1666667 characters of the Normal font serialize to a stored
# OOXML column width of exactly 18.)
$ws.Columns("F").ColumnWidth = 17.1666666666667

# Move the active selection to F7, matching the latest user interaction.
$ws.Range("F7").Select()
